$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "How can I chop onions without crying?"
$ws.Range("A3").Value = "Why add salt to the water when cooking pasta?"
$ws.Range("A4").Value = "How do you properly cook a steak?"
$ws.Range("A5").Value = "How do you peel garlic easily?"
$ws.Range("A6").Value = "How can I keep pasta from sticking to itself?"

$ws.Range("A6").Select()
